# Adds four new experiment result blocks to the "3 trials both models" sheet:
#   - "Type - Logical Augmentation LXMERT new way"      (rows 41-44, cols A:I)
#   - "Type - Contrastive LXMERT new way 0.01"            (rows 46-49, cols A:I)
#   - "Type - Contrastive LXMERT new way 0.05"            (rows 46-49, cols K:S)
#   - "Type - Contrastive LXMERT new way 0.1"             (rows 46-49, cols U:AC)
# Mirrors the existing 5-row block layout used throughout the sheet:
#   header row (type name + column headers), 3 run rows (B/C raw values,
#   D:I summary formulas on the first of the 3 rows), blank spacer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41
$ws.Range('A41').Value = 'Type - Logical Augmentation LXMERT new way'
$ws.Range('B41').Value = 'Test acc.'
$ws.Range('C41').Value = 'Dev acc.'
$ws.Range('D41').Value = 'Min test'
$ws.Range('E41').Value = 'Max test'
$ws.Range('F41').Value = 'Mean test'
$ws.Range('G41').Value = 'Min dev'
$ws.Range('H41').Value = 'Max dev'
$ws.Range('I41').Value = 'Mean dev'

# Row 42
$ws.Range('A42').Value = 'Run#1'
$ws.Range('B42').Value = 77.37
$ws.Range('C42').Value = 74.48
$ws.Range('D42').Formula = '=MIN(B42:B46)'
$ws.Range('E42').Formula = '=MAX(B42:B46)'
$ws.Range('F42').Formula = '=AVERAGE(B42:B46)'
$ws.Range('G42').Formula = '=MIN(C42:C46)'
$ws.Range('H42').Formula = '=MAX(C42:C46)'
$ws.Range('I42').Formula = '=AVERAGE(C42:C46)'

# Row 43
$ws.Range('A43').Value = 'Run#2'
$ws.Range('B43').Value = 78.36
$ws.Range('C43').Value = 76.56
$ws.Range('F43').Formula = '=MAX(F42-D42,E42-F42)'
$ws.Range('I43').Formula = '=MAX(H42-I42,I42-G42)'

# Row 44
$ws.Range('A44').Value = 'Run#3'
$ws.Range('B44').Value = 78.45
$ws.Range('C44').Value = 77.25

# Row 46
$ws.Range('A46').Value = 'Type - Contrastive LXMERT new way 0.01'
$ws.Range('B46').Value = 'Test acc.'
$ws.Range('C46').Value = 'Dev acc.'
$ws.Range('D46').Value = 'Min test'
$ws.Range('E46').Value = 'Max test'
$ws.Range('F46').Value = 'Mean test'
$ws.Range('G46').Value = 'Min dev'
$ws.Range('H46').Value = 'Max dev'
$ws.Range('I46').Value = 'Mean dev'
$ws.Range('K46').Value = 'Type - Contrastive LXMERT new way 0.05'
$ws.Range('L46').Value = 'Test acc.'
$ws.Range('M46').Value = 'Dev acc.'
$ws.Range('N46').Value = 'Min test'
$ws.Range('O46').Value = 'Max test'
$ws.Range('P46').Value = 'Mean test'
$ws.Range('Q46').Value = 'Min dev'
$ws.Range('R46').Value = 'Max dev'
$ws.Range('S46').Value = 'Mean dev'
$ws.Range('U46').Value = 'Type - Contrastive LXMERT new way 0.1'
$ws.Range('V46').Value = 'Test acc.'
$ws.Range('W46').Value = 'Dev acc.'
$ws.Range('X46').Value = 'Min test'
$ws.Range('Y46').Value = 'Max test'
$ws.Range('Z46').Value = 'Mean test'
$ws.Range('AA46').Value = 'Min dev'
$ws.Range('AB46').Value = 'Max dev'
$ws.Range('AC46').Value = 'Mean dev'

# Row 47
$ws.Range('A47').Value = 'Run#1'
$ws.Range('B47').Value = 78.11
$ws.Range('C47').Value = 76.760000000000005
$ws.Range('D47').Formula = '=MIN(B47:B51)'
$ws.Range('E47').Formula = '=MAX(B47:B51)'
$ws.Range('F47').Formula = '=AVERAGE(B47:B51)'
$ws.Range('G47').Formula = '=MIN(C47:C51)'
$ws.Range('H47').Formula = '=MAX(C47:C51)'
$ws.Range('I47').Formula = '=AVERAGE(C47:C51)'
$ws.Range('K47').Value = 'Run#1'
$ws.Range('L47').Value = 77.56
$ws.Range('M47').Value = 77.650000000000006
$ws.Range('N47').Formula = '=MIN(L47:L51)'
$ws.Range('O47').Formula = '=MAX(L47:L51)'
$ws.Range('P47').Formula = '=AVERAGE(L47:L51)'
$ws.Range('Q47').Formula = '=MIN(M47:M51)'
$ws.Range('R47').Formula = '=MAX(M47:M51)'
$ws.Range('S47').Formula = '=AVERAGE(M47:M51)'
$ws.Range('U47').Value = 'Run#1'
$ws.Range('V47').Value = 77.459999999999994
$ws.Range('W47').Value = 76.66
$ws.Range('X47').Formula = '=MIN(V47:V51)'
$ws.Range('Y47').Formula = '=MAX(V47:V51)'
$ws.Range('Z47').Formula = '=AVERAGE(V47:V51)'
$ws.Range('AA47').Formula = '=MIN(W47:W51)'
$ws.Range('AB47').Formula = '=MAX(W47:W51)'
$ws.Range('AC47').Formula = '=AVERAGE(W47:W51)'

# Row 48
$ws.Range('A48').Value = 'Run#2'
$ws.Range('B48').Value = 77.41
$ws.Range('C48').Value = 77.45
$ws.Range('F48').Formula = '=MAX(F47-D47,E47-F47)'
$ws.Range('I48').Formula = '=MAX(H47-I47,I47-G47)'
$ws.Range('K48').Value = 'Run#2'
$ws.Range('L48').Value = 77.66
$ws.Range('M48').Value = 76.36
$ws.Range('P48').Formula = '=MAX(P47-N47,O47-P47)'
$ws.Range('S48').Formula = '=MAX(R47-S47,S47-Q47)'
$ws.Range('U48').Value = 'Run#2'
$ws.Range('V48').Value = 78.599999999999994
$ws.Range('W48').Value = 74.88
$ws.Range('Z48').Formula = '=MAX(Z47-X47,Y47-Z47)'
$ws.Range('AC48').Formula = '=MAX(AB47-AC47,AC47-AA47)'

# Row 49
$ws.Range('A49').Value = 'Run#3'
$ws.Range('B49').Value = 77.459999999999994
$ws.Range('C49').Value = 77.94
$ws.Range('K49').Value = 'Run#3'
$ws.Range('L49').Value = 76.97
$ws.Range('M49').Value = 76.06
$ws.Range('U49').Value = 'Run#3'
$ws.Range('V49').Value = 78.06
$ws.Range('W49').Value = 77.25

# Update the active cell selection to match the new bottom of the sheet.
$ws.Range('I45').Select()
